$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Rewrite the lead-in sentence before the hyperlink in paragraph 3
#    ("Currently consulting with " -> the new "He is co-founder ... for the ")
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Currently consulting with ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Text = "He is co-founder and director of Malstow Geospatial, a consultancy firm offering bespoke consulting and services in the geospatial, geotechnology, maps and location based services fields. This means Gary is currently consulting as Head of APIs for the "

# ---------------------------------------------------------------------------
# 2. Swap the hyperlink target/text from "Lokku" -> "Ordnance Survey"
# ---------------------------------------------------------------------------
$h = $d.Hyperlinks(1)
$h.TextToDisplay = "Ordnance Survey"
$h.Address = "https://www.ordnancesurvey.co.uk"

# ---------------------------------------------------------------------------
# 3. Replace the trailing text after the hyperlink with the new closing
#    clause.
# ---------------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute(" as Geotechnologist in Residence, Gary is helping to advance open geospatial technologies and bring them to new markets.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r2.Text = ", the United Kingdom's national mapping agency."

# ---------------------------------------------------------------------------
# 4. Fix the "APIS as w|ell as" split further down in the bio into a single
#    run of "APIS as well as" (this also relocates the stray _GoBack
#    bookmark that used to live inside that split).
# ---------------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute(" APIS as w", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r3.Text = " APIS as well as "
$r4 = $d.Content
$r4.Find.Execute("ell as ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r4.Text = ""

# ---------------------------------------------------------------------------
# 5. Re-plant the _GoBack bookmark immediately after the "Ordnance Survey"
#    hyperlink (collapsed / empty range), matching its new location.
# ---------------------------------------------------------------------------
$r5 = $d.Content
$r5.Find.Execute(", the United Kingdom's national mapping agency.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmRange = $d.Range($r5.Start, $r5.Start)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
